# Nononcology smoke test data
# Replace the "Economic" / scenario2-3-4 test block (rows 6-14) with the
# simplified "Nononcology" scenario2 block (rows 6-10), and update the
# I5 cell that used to reference the Economic report name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out everything from row 6 through the old bottom of the sheet (row 14)
$ws.Range("A6:J14").ClearContents()

# --- Row 5 / 6 (ExpectedFilenames column I) ---
$ws.Range("I5").Value = "CompleteExcelReport-LIVEHTA Automation-Test_NonOncology_Automation_3-Clinical-2023_"
$ws.Range("I6").Value = "StandardExcelReport-LIVEHTA Automation-Test_NonOncology_Automation_3-Clinical-2023_"

# --- Row 7: scenario2 header data ---
$ws.Range("A7").Value = "scenario2"
$ws.Range("B7").Value = "LIVEHTA Automation - Test_NonOncology_Automation_3"
$ws.Range("C7").Value = "LIVEHTA Automation - Test_NonOncology_Automation_3_radio_button"
$ws.Range("D7").Value = "Clinical"
$ws.Range("E7").Value = "Clinical_radio_button"
$ws.Range("F7").Value = "sub_pop_section1"
$ws.Range("G7").Value = "sub_pop_section1_checkbox"
$ws.Range("H7").Value = "sub_pop_section"

# --- Row 8 ---
# A8 used to be an empty, centre-styled cell (style index 1); it now holds
# a plain, unstyled value, so clear its style before setting the value.
$ws.Range("A8").Style = "Normal"
$ws.Range("A8").Value = "scenario2"
$ws.Range("F8").Value = "intervention_section4"
$ws.Range("G8").Value = "intervention_section4_checkbox"
$ws.Range("H8").Value = "intervention_section"

# --- Row 9 ---
$ws.Range("A9").Value = "scenario2"
$ws.Range("F9").Value = "study_design_section1"
$ws.Range("G9").Value = "study_design_section1_checkbox"
$ws.Range("H9").Value = "study_design_section"

# --- Row 10 ---
$ws.Range("A10").Value = "scenario2"
$ws.Range("F10").Value = "reported_variable_section3"
$ws.Range("G10").Value = "reported_variable_section3_checkbox"
$ws.Range("H10").Value = "reported_variable_section"

# Update sheet view to match the saved selection/scroll position
$ws.Range("I5:I6").Select()
$excel.ActiveWindow.ScrollColumn = 4
